# New crime data collected -- weekly CompStat refresh (76th Precinct).
# Updates the report header (volume/week dates) and the Crime Complaints
# table (rows 16-30) with the newly collected weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
# "Volume 32   Number  42" -> "...43"
$ws.Range("A8").Value = "Volume 32   Number  43"
# "Report Covering the Week  10/13/2025  Through  10/19/2025" -> next week
$ws.Range("C9").Value = "Report Covering the Week  10/20/2025  Through  10/26/2025"

# --- Cells that flip from a text placeholder ("0") to a real number --
# These cells currently hold the shared "0" text string (style carries
# no number format). Give them a numeric format first so the engine
# treats the new value as a genuine number and reuses the existing
# "#,##0" numeric style instead of minting a text style.
foreach ($ref in @("C22","F22","C28","F28")) {
    $ws.Range($ref).Value = 1
    $ws.Range($ref).NumberFormat = "#,##0"
}

# --- Cell that flips from a number back to the text placeholder "0" --
# Write the literal text first (quote-prefix forces text, not a number),
# then pull the formatting back from a cell that already has the plain
# "right/top, no number format" text style so it matches the other
# placeholder cells exactly.
$ws.Range("C25").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C25").PasteSpecial(-4122)

# --- Plain numeric updates --------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 43
$ws.Range("J16").Value = 53
$ws.Range("K16").Value = -18.867924528301
$ws.Range("L16").Value = -23.214285714285
$ws.Range("M16").Value = -50.574712643678
$ws.Range("N16").Value = -86.349206349206
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 114.285714285714
$ws.Range("I17").Value = 89
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = -11.881188118811
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 12.658227848101
$ws.Range("N17").Value = -65.76923076923
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 74
$ws.Range("K18").Value = 20.27027027027
$ws.Range("L18").Value = 17.105263157894
$ws.Range("M18").Value = -1.111111111111
$ws.Range("N18").Value = -74.643874643874
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 12
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 165
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = -2.941176470588
$ws.Range("L19").Value = 24.060150375939
$ws.Range("M19").Value = -12.698412698412
$ws.Range("N19").Value = -16.243654822335
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -85.714285714285
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = 51
$ws.Range("K20").Value = -35.294117647058
$ws.Range("L20").Value = -35.294117647058
$ws.Range("M20").Value = -37.735849056603
$ws.Range("N20").Value = -90.207715133531
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 41
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 36.666666666666
$ws.Range("I21").Value = 423
$ws.Range("J21").Value = 454
$ws.Range("K21").Value = -6.828193832599
$ws.Range("L21").Value = 2.919708029197
$ws.Range("M21").Value = -15.568862275449
$ws.Range("N21").Value = -71.341463414634
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 92
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 2.222222222222
$ws.Range("L23").Value = 3.370786516853
$ws.Range("M23").Value = 35.294117647058
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -70.588235294117
$ws.Range("F24").Value = 35
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = -37.5
$ws.Range("I24").Value = 391
$ws.Range("J24").Value = 447
$ws.Range("K24").Value = -12.527964205816
$ws.Range("L24").Value = -21.956087824351
$ws.Range("M24").Value = 4.545454545454
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -26.086956521739
$ws.Range("I25").Value = 148
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = -33.031674208144
$ws.Range("L25").Value = -42.412451361867
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("I26").Value = 161
$ws.Range("J26").Value = 135
$ws.Range("K26").Value = 19.259259259259
$ws.Range("L26").Value = 24.806201550387
$ws.Range("M26").Value = -33.195020746888
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 9.090909090909
$ws.Range("N29").Value = -94.230769230769
$ws.Range("N30").Value = -95.121951219512

Write-Output "Applied weekly CompStat update"
